# Update the "Förändrad" date column (C) for rows 2-6 from 2023-10-13 (45212)
# to 2023-10-22 (45221), matching the automatic update of the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$newValue = 45221

for ($row = 2; $row -le 6; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45212) {
        $cell.Value = $newValue
    }
}
